$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price (column D) cells: these hold numeric-looking values that must stay
# stored as TEXT (matching the workbook's inlineStr convention). Marking the
# cell's NumberFormat as Text ("@") before assigning the value prevents
# Excel's COM layer from auto-coercing the string into a real number; the
# Style reset afterwards clears the now-unneeded direct formatting so the
# cell keeps its original (unstyled) appearance.
$priceUpdates = @{
    "D2"  = "244.93"
    "D4"  = "5.412"
    "D5"  = "0.05983"
    "D6"  = "3.390"
    "D8"  = "0.9270"
    "D9"  = "0.1423"
    "D10" = "0.07441"
    "D11" = "0.03367"
    "D12" = "0.03037"
    "D13" = "0.09360"
    "D14" = "3.937"
    "D15" = "0.001601"
    "D16" = "0.04828"
    "D18" = "0.005661"
    "D20" = "0.0009822"
    "D21" = "0.00007703"
    "D22" = "3.659"
    "D40" = "0.03944"
    "D41" = "0.006212"
}

foreach ($addr in $priceUpdates.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $priceUpdates[$addr]
    $cell.Style = "Normal"
}

# Volume(1h) label (column E) cells: plain text swaps, no numeric coercion
# risk so a direct Value assignment is fine.
$ws.Range("E17").Value = "16OneONE"
$ws.Range("E41").Value = "40KickTokenKICKBestin24h"
$ws.Range("E44").Value = "43LocalTradersLCT"
$ws.Range("E49").Value = "48BOLOBOLOWorstin24h"
